$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (38 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 891
$ws.Range("I12").Value = 749.375
$ws.Range("K12").Value = 749.375
$ws.Range("M12").Value = -579.375
$ws.Range("H40").Value = 1152679.2
$ws.Range("I40").Value = 2980.8333
$ws.Range("J40").Value = 1964231.1
$ws.Range("K40").Value = 2980.8333
$ws.Range("L40").Value = 1964231.1
$ws.Range("M40").Value = -2805.8333
$ws.Range("N40").Value = -1964581.1
$ws.Range("H86").Value = 41714980
$ws.Range("J86").Value = 8001733
$ws.Range("L86").Value = 8001733
$ws.Range("N86").Value = -8003979
$ws.Range("H89").Value = 41714980
$ws.Range("J89").Value = 8001733
$ws.Range("L89").Value = 40008665
$ws.Range("N89").Value = -40019897
$ws.Range("H113").Value = 70059980
$ws.Range("I113").Value = 12348597
$ws.Range("J113").Value = 107160150
$ws.Range("K113").Value = 12348597
$ws.Range("L113").Value = 107160150
$ws.Range("M113").Value = -12345343
$ws.Range("N113").Value = -107166658
$ws.Range("H115").Value = 715.63635
$ws.Range("I115").Value = 687.3
$ws.Range("K115").Value = 2061.9
$ws.Range("M115").Value = -494.8999999999996
$ws.Range("H132").Value = 1299.9868
$ws.Range("I132").Value = 1274.7123
$ws.Range("K132").Value = 3824.1369
$ws.Range("M132").Value = -1294.1369
$ws.Range("H141").Value = 6398.75
$ws.Range("I141").Value = 6398.75
$ws.Range("K141").Value = 19196.25
$ws.Range("M141").Value = -14016.25

# --- Sheet: ARM (15 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2451829.5
$ws.Range("I97").Value = 595.1
$ws.Range("K97").Value = 595.1
$ws.Range("M97").Value = -99.10000000000002
$ws.Range("H110").Value = 8333799
$ws.Range("I110").Value = 470.5484
$ws.Range("K110").Value = 470.5484
$ws.Range("M110").Value = 1574.4516
$ws.Range("H132").Value = 5777.9556
$ws.Range("I132").Value = 3158.1428
$ws.Range("J132").Value = 10092.941
$ws.Range("K132").Value = 9474.428400000001
$ws.Range("L132").Value = 30278.823
$ws.Range("M132").Value = -6944.428400000001
$ws.Range("N132").Value = -35338.823

# --- Sheet: BSM (4 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7761.3335
$ws.Range("I134").Value = 2491.8
$ws.Range("K134").Value = 7475.400000000001
$ws.Range("M134").Value = -4940.400000000001

# --- Sheet: CRP (11 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4603.3096
$ws.Range("I132").Value = 2395.4583
$ws.Range("K132").Value = 7186.374899999999
$ws.Range("M132").Value = -4656.374899999999
$ws.Range("H134").Value = 6870.706
$ws.Range("I134").Value = 3059.75
$ws.Range("J134").Value = 8949.409
$ws.Range("K134").Value = 9179.25
$ws.Range("L134").Value = 26848.227
$ws.Range("M134").Value = -6644.25
$ws.Range("N134").Value = -31918.227

# --- Sheet: CUL (34 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3497.111
$ws.Range("J5").Value = 6255.222
$ws.Range("L5").Value = 18765.666
$ws.Range("N5").Value = -18989.666
$ws.Range("H64").Value = 3943.9167
$ws.Range("I64").Value = 730
$ws.Range("K64").Value = 2190
$ws.Range("M64").Value = -1920
$ws.Range("H67").Value = 3943.9167
$ws.Range("I67").Value = 730
$ws.Range("K67").Value = 2190
$ws.Range("M67").Value = -1254
$ws.Range("H113").Value = 6644.231
$ws.Range("J113").Value = 7114.5835
$ws.Range("L113").Value = 21343.7505
$ws.Range("N113").Value = -25683.7505
$ws.Range("H135").Value = 3497.111
$ws.Range("J135").Value = 6255.222
$ws.Range("L135").Value = 56296.998
$ws.Range("N135").Value = -61366.998
$ws.Range("H137").Value = 85549.46000000001
$ws.Range("I137").Value = 169221.33
$ws.Range("J137").Value = 57658.832
$ws.Range("K137").Value = 507663.99
$ws.Range("L137").Value = 172976.496
$ws.Range("M137").Value = -502563.99
$ws.Range("N137").Value = -183176.496
$ws.Range("H139").Value = 95988.55
$ws.Range("I139").Value = 146726
$ws.Range("J139").Value = 7198
$ws.Range("K139").Value = 440178
$ws.Range("L139").Value = 21594
$ws.Range("M139").Value = -435038
$ws.Range("N139").Value = -31874

# --- Sheet: GSM (22 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4880.6665
$ws.Range("I80").Value = 2847.8333
$ws.Range("J80").Value = 6913.5
$ws.Range("K80").Value = 2847.8333
$ws.Range("L80").Value = 6913.5
$ws.Range("M80").Value = -1849.8333
$ws.Range("N80").Value = -8909.5
$ws.Range("H83").Value = 4880.6665
$ws.Range("I83").Value = 2847.8333
$ws.Range("J83").Value = 6913.5
$ws.Range("K83").Value = 14239.1665
$ws.Range("L83").Value = 34567.5
$ws.Range("M83").Value = -9247.166499999999
$ws.Range("N83").Value = -44551.5
$ws.Range("H97").Value = 928.6957
$ws.Range("J97").Value = 1056.875
$ws.Range("L97").Value = 1056.875
$ws.Range("N97").Value = -2048.875
$ws.Range("H122").Value = 4846065.5
$ws.Range("J122").Value = 6182.6665
$ws.Range("L122").Value = 18547.9995
$ws.Range("N122").Value = -23447.9995

# --- Sheet: LTW (27 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6999.8
$ws.Range("I7").Value = 4499.5
$ws.Range("K7").Value = 4499.5
$ws.Range("M7").Value = -4387.5
$ws.Range("H55").Value = 390.86365
$ws.Range("I55").Value = 21.88889
$ws.Range("J55").Value = 646.3077
$ws.Range("K55").Value = 21.88889
$ws.Range("L55").Value = 646.3077
$ws.Range("M55").Value = 151.11111
$ws.Range("N55").Value = -992.3077
$ws.Range("H61").Value = 4903.657
$ws.Range("I61").Value = 3779.3333
$ws.Range("K61").Value = 3779.3333
$ws.Range("M61").Value = -3577.3333
$ws.Range("H113").Value = 4903.657
$ws.Range("I113").Value = 3779.3333
$ws.Range("K113").Value = 3779.3333
$ws.Range("M113").Value = -1609.3333
$ws.Range("H126").Value = 6999.8
$ws.Range("I126").Value = 4499.5
$ws.Range("K126").Value = 13498.5
$ws.Range("M126").Value = -11028.5
$ws.Range("H136").Value = 14529.277
$ws.Range("I136").Value = 3539.4285
$ws.Range("K136").Value = 10618.2855
$ws.Range("M136").Value = -8068.2855

# --- Sheet: WVR (22 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 71951.57000000001
$ws.Range("I62").Value = 71951.57000000001
$ws.Range("K62").Value = 71951.57000000001
$ws.Range("M62").Value = -71327.57000000001
$ws.Range("H65").Value = 71951.57000000001
$ws.Range("I65").Value = 71951.57000000001
$ws.Range("K65").Value = 359757.85
$ws.Range("M65").Value = -356637.85
$ws.Range("I81").Value = 1041331.3
$ws.Range("K81").Value = 2082662.6
$ws.Range("M81").Value = -2081601.6
$ws.Range("I84").Value = 1041331.3
$ws.Range("K84").Value = 10413313
$ws.Range("M84").Value = -10408009
$ws.Range("H126").Value = 71432760
$ws.Range("I126").Value = 111114470
$ws.Range("K126").Value = 333343410
$ws.Range("M126").Value = -333340940
$ws.Range("H136").Value = 337179.5
$ws.Range("I136").Value = 1131.9445
$ws.Range("K136").Value = 3395.8335
$ws.Range("M136").Value = -845.8335000000002
